$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.182.62"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").Value = "1.643.61"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'216.88"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").Value = "'19.76"
$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "1.871.45"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("D14").Value = "1.650.39"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").Value = "'0.544"
$ws.Range("E15").Value = "  -3.25%  "

$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").Value = "'63.22"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "26.178.82"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").Value = "'194.96"
$ws.Range("E21").Value = "  +1.30%  "

$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").Value = "'142.62"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("E27").Value = "  +1.09%  "

$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "'15.60"
$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Value = "'3.24"
$ws.Range("E33").Value = "  +0.54%  "

$ws.Range("E34").Value = "  +1.60%  "

$ws.Range("E35").Value = "  +1.45%  "

$ws.Range("D36").Value = "'0.911"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("D37").Value = "1.132.21"
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "'0.552"
$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("E39").Value = "  -1.17%  "

$ws.Range("E40").Value = "  +1.15%  "

$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("D42").Value = "'100.40"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("D44").Value = "'0.796"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").Value = "1.780.55"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").Value = "'56.83"
$ws.Range("E47").Value = "  +2.69%  "

$ws.Range("E48").Value = "  +3.41%  "

$ws.Range("E49").Value = "  +2.89%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.69"
$ws.Range("E50").Value = "  +3.30%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.417"
$ws.Range("E51").Value = "  +0.08%  "
